$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric must be forced to stay as text,
# matching the source data which stores them as literal strings
# (e.g. '1.001', '28.370.51') rather than numbers.
$textForceCells = @('D4', 'D5', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D14', 'D15', 'D16', 'D18', 'D19', 'D20', 'D21', 'D23', 'D24', 'D25', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.370.51'
$ws.Range("E2").Value = '  +3.42%  '
$ws.Range("D3").Value = '1.869.20'
$ws.Range("E3").Value = '  +1.83%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '339.26'
$ws.Range("E5").Value = '  +1.89%  '
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("D7").Value = '0.4692'
$ws.Range("E7").Value = '  +1.54%  '
$ws.Range("D8").Value = '0.3973'
$ws.Range("E8").Value = '  +4.13%  '
$ws.Range("D9").Value = '47.51'
$ws.Range("E9").Value = '  +2.12%  '
$ws.Range("D10").Value = '0.08020'
$ws.Range("E10").Value = '  +1.68%  '
$ws.Range("D11").Value = '0.9993'
$ws.Range("E11").Value = '  +2.76%  '
$ws.Range("D12").Value = '21.95'
$ws.Range("E12").Value = '  +4.08%  '
$ws.Range("D13").Value = '1.871.62'
$ws.Range("E13").Value = '  +1.26%  '
$ws.Range("D14").Value = '6.022'
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("D15").Value = '7.247'
$ws.Range("E15").Value = '  +3.08%  '
$ws.Range("D16").Value = '91.21'
$ws.Range("E16").Value = '  +3.85%  '
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").Value = '0.00001043'
$ws.Range("E18").Value = '  +1.51%  '
$ws.Range("D19").Value = '0.06613'
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").Value = '17.55'
$ws.Range("E20").Value = '  +3.53%  '
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("D22").Value = '28.374.93'
$ws.Range("E22").Value = '  +3.38%  '
$ws.Range("D23").Value = '5.474'
$ws.Range("E23").Value = '  +2.28%  '
$ws.Range("D24").Value = '11.06'
$ws.Range("E24").Value = '  +1.95%  '
$ws.Range("D25").Value = '2.258'
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("D26").Value = '2.085.22'
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("D27").Value = '161.16'
$ws.Range("E27").Value = '  +2.41%  '
$ws.Range("D28").Value = '19.80'
$ws.Range("E28").Value = '  +2.32%  '
$ws.Range("D29").Value = '2.126'
$ws.Range("E29").Value = '  +2.93%  '
$ws.Range("D30").Value = '5.502'
$ws.Range("E30").Value = '  +3.44%  '
$ws.Range("D31").Value = '120.29'
$ws.Range("E31").Value = '  +1.14%  '
$ws.Range("D32").Value = '0.9681'
$ws.Range("E32").Value = '  +1.57%  '
$ws.Range("D33").Value = '0.09494'
$ws.Range("E33").Value = '  +2.11%  '
$ws.Range("D34").Value = '3.590'
$ws.Range("E34").Value = '  +0.69%  '
$ws.Range("D35").Value = '5.351'
$ws.Range("E35").Value = '  +2.11%  '
$ws.Range("D36").Value = '1.371'
$ws.Range("E36").Value = '  +4.07%  '
$ws.Range("D37").Value = '0.06096'
$ws.Range("E37").Value = '  +2.69%  '
$ws.Range("D38").Value = '0.02249'
$ws.Range("E38").Value = '  +2.69%  '
$ws.Range("D39").Value = '8.361'
$ws.Range("E39").Value = '  +3.57%  '
$ws.Range("D40").Value = '1.183'
$ws.Range("E40").Value = '  +2.12%  '
$ws.Range("D41").Value = '0.5945'
$ws.Range("E41").Value = '  +2.63%  '
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("D43").Value = '0.1874'
$ws.Range("E43").Value = '  +1.81%  '
$ws.Range("D44").Value = '10.36'
$ws.Range("E44").Value = '  +3.48%  '
$ws.Range("D45").Value = '1.291'
$ws.Range("E45").Value = '  +2.46%  '
$ws.Range("D46").Value = '0.5589'
$ws.Range("E46").Value = '  +1.92%  '
$ws.Range("D47").Value = '12.20'
$ws.Range("E47").Value = '  +1.25%  '
$ws.Range("D48").Value = '1.959'
$ws.Range("E48").Value = '  +4.87%  '
$ws.Range("D49").Value = '0.06868'
$ws.Range("E49").Value = '  +3.14%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '2.055'
$ws.Range("E50").Value = '  +16.45%  '
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").Value = '111.75'
$ws.Range("E51").Value = '  +1.69%  '
